$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Sheet1: the active-sheet selection changes to A1:B1 (and it stops being tabSelected
# once Sheet2 becomes active, further below)
$sheet1.Range("A1:B1").Select()

# Insert the new worksheet right after Sheet1, named Sheet2
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$newSheet.Name = "Sheet2"

# B/C content for each data row (Column A is the constant "HELP" label)
$rows = @(
    @('Alternativas de transición a post escuela', 'Ahora viene lo mejor! Herramientas para una transición positiva'),
    @('Alternativas de transición a post escuela', 'Mapa de Ruta al éxito postsecundario'),
    @('Alternativas de transición a post escuela', 'Un camino de triunfo para mi hijo'),
    @('Apoyo a la adquisición de idiomas  ', 'Konesans se pouvwa (Conocimiento es poder) Lo que el saber otras lenguas puede hacer por tu hijo/a'),
    @('Apoyo a la transición escolar y ambientes nuevos', 'Las transiciones… el pan nuestro de cada día'),
    @('Apoyo a la transición escolar y ambientes nuevos', 'Preparación para la escuela y las transiciones'),
    @('Apoyo a la transición escolar y ambientes nuevos', 'La transición escolar y la maravilla de descubrir nuevas posibilidades'),
    @('Cómo desarrollar buenos hábitos de estudio', 'Cómo preparamos a nuestros hijos a tener buenos hábitos de estudio?'),
    @('Cómo desarrollar buenos hábitos de estudio', 'Eb= Em. Estudiar bien es igual a ejecutar mejor: el desarrollo de hábitos de estudio'),
    @('Cómo desarrollar buenos hábitos de estudio', 'Plan de estudio de mi hijo'),
    @('Cómo desarrollar buenos hábitos de estudio. Apoyo a la adquisición de idiomas', 'Explorando la tecnología de videojuegos controlados por la mente'),
    @('Cómo interpretar exámenes y el informe de rendimiento académico. Cómo monitorear el progreso del estudiante', 'Entendiendo tu esfuerzo: lo que me hablas a través de tus calificaciones'),
    @('Cómo monitorear el progreso del estudiante ', 'Siguiendo el camino correcto: mis oportunidades para monitorear el progreso de mi hijo/a'),
    @('Cómo promover la lectura en el hogar', 'Círculos literarios para padres'),
    @('Cómo promover la lectura en el hogar', 'La magia de la escritura'),
    @('Cómo promover la lectura en el hogar', 'Maneras sencillas de ayudar a desarrollar la lectura en mis hijos desde su nacimiento'),
    @('Cómo promover la lectura en el hogar', 'Si leemos, entendemos: el arte de la lectura en el ambiente hogareño'),
    @('Cómo promover la lectura en el hogar. Apoyo a la adquisición de idiomas ', 'Mi nene lee…pero no entiende lo que lee. La importancia de la comprensión lectora en el aprovechamiento académico de mi hijo, de mi nieto.'),
    @('Comunicación efectiva', 'Mis padres discuten constantemente: pensamientos de un hijo (a)'),
    @('Comunicación efectiva. Promoción de valores en el hogar ', 'El juego simbólico como campo de aprendizaje – Enfoque Sistémico'),
    @('Comunicación efectiva. Técnicas para el manejo de la crianza ', 'El arte, el cerebro y nuestros pensamientos”'),
    @('Comunicación Efectiva. Apoyo a la transición', 'La experiencia de comunicarte hábilmente con tu hijo y su relación en su desempeño escolar'),
    @('Comunicación efectiva. Cómo monitorear el progreso del estudiante', 'Escuchando lo que no se dice'),
    @('Comunicación efectiva. Promoción de valores en el hogar', 'Comunicación efectiva entre miembros del “equipo invencible”'),
    @('Comunicación efectiva. Técnicas de disciplina positiva. Promoción de valores.', 'Grupo de apoyo- Compartiendo experiencias '),
    @('Curso de contenido de idioma ', 'Desde aquí hasta el infinito: el panorama universal ante el dominio de idiomas'),
    @('Establecer y mantener redes de apoyo en la comunidad', 'Aplicación práctica de la Pedagogía Sistémica en el escenario familiar'),
    @('Establecer y mantener redes de apoyo en la comunidad', 'El padre como gestor de la escuela para mejorar el aprovechamiento de su hijo.'),
    @('Establecer y mantener redes de apoyo en la comunidad', 'La biblioteca pública de mi pueblo… cómo me puede ayudar en la educación de mi hijo.'),
    @('Establecer y mantener redes de apoyo en la comunidad', 'Mi comunidad…mi gran aliado'),
    @('Establecer y mantener redes de apoyo en la comunidad', 'Explicando la Carta Circular de participación de padres'),
    @('Establecer y mantener redes de apoyo en la comunidad. Liderazgo. Cómo monitorear el progreso del estudiante.', 'A decir presente y a participar!'),
    @('Factores de riesgo y prevención', 'Cómo trabajar para tener hijos menos violentos y agresivos'),
    @('Factores de riesgo y prevención', 'Conociendo los rasgos de personalidad de mi hijo (a) para optimizar la comunicación'),
    @('Factores de riesgo y prevención', 'La inteligencia emocional y sus implicaciones en el éxito de vida de mi hijo'),
    @('Factores de riesgo y prevención', 'Las redes sociales en el entorno de mi hijo'),
    @('Factores de riesgo y prevención', 'Prevención de Violencia Doméstica'),
    @('Factores de riesgo y prevención', 'Resiliencia ¿Cómo lograr la excelencia en el lugar de trabajo?'),
    @('Factores de riesgo y prevención ', 'Protegiendo mi semilla'),
    @('Factores de riesgo y prevención  ', 'Cómo se lo explico a mi hijo?'),
    @('Factores de riesgo y prevención. Comunicación efectiva.', 'Está bien!… “Si no se siente bien”'),
    @('Factores de riesgo y prevención. Comunicación efectiva.', 'Los sentimientos de mi hijo'),
    @('Factores de riesgo y prevención. Leyes que protegen la niñez (Educación Especial, “cyberbullying”, etc.).', 'Acoso escolar, “Bullying”'),
    @('Leyes que protegen la niñez (Educación Especial, “cyberbullying”, etc.).', 'El derecho a mi favor: leyes que protegen a mis hijos '),
    @('Leyes que protegen la niñez (Educación Especial, “cyberbullying”, etc.).', 'El derecho penal como respuesta a la actividad delictiva '),
    @('Liderazgo', 'Promoviendo la responsabilidad y autonomía de mi hijo '),
    @('Liderazgo ', 'El líder que hay en ti: alcanzando la plenitud mediante el ejercicio del liderazgo'),
    @('Mediación de conflictos. Comunicación efectiva', 'Conflictos: Oportunidades para conocer y crecer.'),
    @('Mediación de conflictos. Comunicación efectiva', 'Ese conflicto… tema superado. Cómo trabajar adecuadamente las situaciones conflictivas'),
    @('Mediación de conflictos. Promoción de valores en el hogar. Factores de riesgo y prevención ', 'Manejo adecuado de conflictos'),
    @('Navegando el  sistema  educativo. Cómo monitorear el progreso del estudiante.', 'Sistema educativo: un mar que debo navegar'),
    @('Navegando el sistema educativo', 'Navegando en la página del Departamento de Educación'),
    @('Navegando el sistema educativo. Cómo monitorear el progreso del estudiante.', 'Sistema de Información Estudiantil (SIE)'),
    @('Preparación. Motivación en la participación de las pruebas del estado. Cómo desarrollar buenos hábitos de estudio.', 'Tu logro es el mío! Preparándote para dar el máximo en las pruebas'),
    @('Promoción de valores en el hogar', 'Autogestión: Bisutería, Diseño y Creación de “Choker”'),
    @('Promoción de valores en el hogar', 'Biografía de mi corazón!'),
    @('Promoción de valores en el hogar', 'Las competencias educativas a través de una mirada sistémica'),
    @('Promoción de valores en el hogar', 'Manejo de la sala de clases a través de valores'),
    @('Promoción de valores en el hogar', 'Presupuesto familiar:  equilibrio entre el dar y el tomar'),
    @('Promoción de valores en el hogar. Comunicación efectiva.', 'Presente en su Presente”'),
    @('Promoción de valores en el hogar. Comunicación efectiva.', 'Mi sistema de valores: la transmisión de valores positivos en el hogar'),
    @('Promoción de valores en el hogar. Comunicación efectiva.', 'Respetando y aceptando la diversidad'),
    @('Promoción de valores en el hogar. Comunicación efectiva.', 'Valores:  la llave del éxito en la vida'),
    @('Técnica para el manejo de la crianza en las etapas de desarrollo', 'La familia como primera escuela de mi hijo '),
    @('Técnicas de disciplina positiva', 'Cuáles son las necesidades reales de mis hijos y cómo saberlo?    '),
    @('Técnicas de disciplina positiva. Comunicación efectiva.', 'El arte de construir con amor'),
    @('Técnicas de disciplina positiva. Comunicación efectiva. Promoción de valores.', 'Manejo de emociones desde la neurociencia'),
    @('Técnicas de disciplina positiva. Promoción de valores en el hogar.', 'Del Padre que soy al Padre que quiero ser.'),
    @('Técnicas de disciplina positiva. Promoción de valores en el hogar. Comunicación efectiva.', 'Disciplina con amor'),
    @('Técnicas para el manejo de crianza en las etapas de desarrollo', 'Te entiendo! Es cuestión del desarrollo….'),
    @('Técnicas para el manejo de crianza en las etapas de desarrollo. Comunicación efectiva.', 'La Adolescencia”'),
    @('Técnicas para el manejo de crianza en las etapas de desarrollo. Comunicación efectiva.', 'La crianza compartida, una alternativa necesaria para la educación de nuestros hijos'),
    @('Técnicas para el manejo de crianza en las etapas de desarrollo. Técnicas para promover la salud, la nutrición y una vida activa de los (as) estudiantes', 'Creciendo y aprendiendo en familia!'),
    @('Técnicas para el manejo de la crianza en las etapas de desarrollo', 'Soy abuelo criando nietos: “necesito herramientas para ayudarlos a estudiar”'),
    @('Técnicas para promover la salud, la nutrición y una vida activa de los (as) estudiantes', 'Actívate, muévete y disfruta de lo saludable previniendo el sobrepeso y obesidad'),
    @('Técnicas para promover la salud, la nutrición y una vida activa de los (as) estudiantes', 'Alimentación sana, mente sana”'),
    @('Técnicas para promover la salud, la nutrición y una vida activa de los (as) estudiantes', 'Lee las etiquetas nutricionales y aprende a seleccionar meriendas saludables'),
    @('Técnicas para promover la salud, la nutrición y una vida activa de los (as) estudiantes', 'Promoviendo una excelente nutrición: “Mi Plato para un Puerto Rico Saludable”'),
    @('Técnicas para promover la salud, la nutrición y una vida activa de los estudiantes', 'La salud del estudiante, trabajo de todos '),
    @('Técnicas para promover la salud, la nutrición y una vida activa de los estudiantes', 'La salud, la buena alimentación y la vida activa como elementos'),
    @('Temas de cultura puertorriqueña ', 'Cultura y esencia: el arte de ser puertorriqueño ')
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $newSheet.Cells.Item($r, 1).Value = 'HELP'
    $newSheet.Cells.Item($r, 2).Value = $rows[$i][0]
    $newSheet.Cells.Item($r, 3).Value = $rows[$i][1]
}

# Column B (rows 1-42) carries a wrap-text style; rows 43-81 stay plain
$newSheet.Range("B1:B42").WrapText = $true

# Rows whose wrapped text needs two lines get an explicit row height
foreach ($r in 11, 12, 25, 32) {
    $newSheet.Rows.Item($r).RowHeight = 30
}

# Column widths (best-fit-like, matching the authored layout)
$newSheet.Columns.Item(2).ColumnWidth = 134.7109375
$newSheet.Columns.Item(3).ColumnWidth = 93.28515625

# Make Sheet2 the active sheet/tab and set its selection + scroll position
$newSheet.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 55
$newSheet.Range("A82").Select()

